$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.04314941262403053
$ws.Cells.Item(2, 8).Value = 2.059558138978731
$ws.Cells.Item(2, 9).Value = 60.79545561192806
$ws.Cells.Item(3, 7).Value = 0.06514116215503429
$ws.Cells.Item(3, 8).Value = 33.75912743349041
$ws.Cells.Item(4, 7).Value = 0.006296131661288136
$ws.Cells.Item(4, 8).Value = 230.6973246056239
$ws.Cells.Item(5, 7).Value = 0.01959938850275934
$ws.Cells.Item(5, 8).Value = 416.2691463811145
$ws.Cells.Item(6, 7).Value = 0.02293909489694262
$ws.Cells.Item(6, 8).Value = -33.83509331140302
$ws.Cells.Item(7, 7).Value = 0.0162694607159024
$ws.Cells.Item(7, 8).Value = -69.41195153833343
$ws.Cells.Item(8, 7).Value = -0.01592594257915191
$ws.Cells.Item(8, 8).Value = 15.3810480266481
$ws.Cells.Item(9, 7).Value = -0.001841710952001857
$ws.Cells.Item(9, 8).Value = 91.44500614069423
$ws.Cells.Item(10, 7).Value = -0.06858393625721131
$ws.Cells.Item(10, 8).Value = 5.663470353926972
$ws.Cells.Item(11, 7).Value = -0.09141119550083635
$ws.Cells.Item(11, 8).Value = 0.6519283113263623
$ws.Cells.Item(12, 7).Value = -0.220474027838467
$ws.Cells.Item(12, 8).Value = 9.809842833187542
$ws.Cells.Item(13, 7).Value = -0.3177153950892906
$ws.Cells.Item(13, 8).Value = -15.61139381129434
$ws.Cells.Item(14, 7).Value = -0.03635396749877619
$ws.Cells.Item(14, 8).Value = 2.003456426700552
$ws.Cells.Item(15, 7).Value = 0.05760701532592673
$ws.Cells.Item(15, 8).Value = 265.6737512594899
$ws.Cells.Item(16, 7).Value = 0.1329894552070179
$ws.Cells.Item(16, 8).Value = 6.129900937740514
$ws.Cells.Item(17, 7).Value = 0.160429742634934
$ws.Cells.Item(17, 8).Value = 14.38800768744737
$ws.Cells.Item(18, 7).Value = 0.1320510714258215
$ws.Cells.Item(18, 8).Value = 5.867223279797911
$ws.Cells.Item(19, 7).Value = 0.1150871047828415
$ws.Cells.Item(19, 8).Value = -13.60724809343366
$ws.Cells.Item(20, 7).Value = 0.04112336532738485
$ws.Cells.Item(20, 8).Value = 19.7670967267929
$ws.Cells.Item(21, 7).Value = 0.05791685721070196
$ws.Cells.Item(21, 8).Value = -0.2119326227162358
$ws.Cells.Item(22, 7).Value = -0.07052924021946347
$ws.Cells.Item(22, 8).Value = 11.67060633660727
$ws.Cells.Item(23, 7).Value = -0.1020486548571485
$ws.Cells.Item(23, 8).Value = -63.15663044877915
$ws.Cells.Item(24, 7).Value = 0.1284528633293741
$ws.Cells.Item(24, 8).Value = 8.752822399611913
$ws.Cells.Item(25, 7).Value = 0.141305506620729
$ws.Cells.Item(25, 8).Value = 11.99577209091021
$ws.Cells.Item(26, 7).Value = 0.04206632653179394
$ws.Cells.Item(26, 8).Value = -15.36649976742151
$ws.Cells.Item(27, 7).Value = 0.09629888295597235
$ws.Cells.Item(27, 8).Value = 11.10139985766141
$ws.Cells.Item(28, 7).Value = -0.07835283892605897
$ws.Cells.Item(28, 8).Value = -23.20784031819758
$ws.Cells.Item(29, 7).Value = -0.08629659800508381
$ws.Cells.Item(29, 8).Value = -21.25039325914446
$ws.Cells.Item(30, 7).Value = 0.06578618681454668
$ws.Cells.Item(30, 8).Value = 3.261844485171043
$ws.Cells.Item(31, 7).Value = 0.05992030089595504
$ws.Cells.Item(31, 8).Value = -1.089724896272676
$ws.Cells.Item(32, 7).Value = 0.07969497738109868
$ws.Cells.Item(32, 8).Value = -18.89706286252383
$ws.Cells.Item(33, 7).Value = 0.09048705000702399
$ws.Cells.Item(33, 8).Value = 9.968249700864444
$ws.Cells.Item(34, 7).Value = -0.002060927227598771
$ws.Cells.Item(34, 8).Value = -107.9098451164058
$ws.Cells.Item(35, 7).Value = -0.01051318589315675
$ws.Cells.Item(35, 8).Value = 6.216915481556747
$ws.Cells.Item(36, 7).Value = 0.01131872201367714
$ws.Cells.Item(36, 8).Value = 2003.04219414925
$ws.Cells.Item(37, 7).Value = 0.00683811696603353
$ws.Cells.Item(37, 8).Value = 154.4684774297722
$ws.Cells.Item(38, 7).Value = 0.1199671450287729
$ws.Cells.Item(38, 8).Value = 11.84989875383528
$ws.Cells.Item(39, 7).Value = 0.1182195955492491
$ws.Cells.Item(39, 8).Value = 38.0057859818868
$ws.Cells.Item(40, 7).Value = 0.0312007359582262
$ws.Cells.Item(40, 8).Value = 950.4388319479764
$ws.Cells.Item(41, 7).Value = 0.02463498740924376
$ws.Cells.Item(41, 8).Value = 64.2688800505002
$ws.Cells.Item(42, 7).Value = 0.1083006493956065
$ws.Cells.Item(42, 8).Value = 7.29728509410414
$ws.Cells.Item(43, 7).Value = 0.105973908647577
$ws.Cells.Item(43, 8).Value = -11.79462579723324
$ws.Cells.Item(44, 7).Value = 0.01785244111131554
$ws.Cells.Item(44, 8).Value = -49.97576867846407
$ws.Cells.Item(45, 7).Value = 0.03160702251619605
$ws.Cells.Item(45, 8).Value = 93.07688605295411
$ws.Cells.Item(46, 7).Value = 0.06864599986596653
$ws.Cells.Item(46, 8).Value = 89.42826994281801
$ws.Cells.Item(47, 7).Value = 0.06635106289238056
$ws.Cells.Item(47, 8).Value = 31.54371534668483
$ws.Cells.Item(48, 7).Value = 0.05843369710995205
$ws.Cells.Item(48, 8).Value = 36.58699962621249
$ws.Cells.Item(49, 7).Value = 0.05726030159442948
$ws.Cells.Item(49, 8).Value = -17.58273558614225
$ws.Cells.Item(50, 7).Value = -0.007451780070311671
$ws.Cells.Item(50, 8).Value = -143.1419365323532
$ws.Cells.Item(51, 7).Value = 0.01750886359214551
$ws.Cells.Item(51, 8).Value = -10.07127147197189
$ws.Cells.Item(52, 7).Value = -0.117805289271679
$ws.Cells.Item(52, 8).Value = -13.79904464555063
$ws.Cells.Item(53, 7).Value = -0.09181723855864715
$ws.Cells.Item(53, 8).Value = 0.5822853592308389
$ws.Cells.Item(54, 7).Value = 0.08505192340343821
$ws.Cells.Item(54, 8).Value = 16.31404300138482
$ws.Cells.Item(55, 7).Value = 0.06874136818618598
$ws.Cells.Item(55, 8).Value = 10.95969203130345
$ws.Cells.Item(56, 7).Value = 0.03673717257822359
$ws.Cells.Item(56, 8).Value = 4.995704092452656
$ws.Cells.Item(57, 7).Value = 0.02351307982883641
$ws.Cells.Item(57, 8).Value = 307.2574285266576
$ws.Cells.Item(58, 7).Value = 0.04896086991176206
$ws.Cells.Item(58, 8).Value = 95.76042281042598
$ws.Cells.Item(59, 7).Value = 0.02174369270121645
$ws.Cells.Item(59, 8).Value = -8.17164432256857
$ws.Cells.Item(60, 7).Value = 0.01982788087313853
$ws.Cells.Item(60, 8).Value = -38.88330335122644
$ws.Cells.Item(61, 7).Value = 0.04678051883948923
$ws.Cells.Item(61, 8).Value = 269.5729505611063
$ws.Cells.Item(62, 7).Value = 0.0534043052293935
$ws.Cells.Item(62, 8).Value = -11.5275197163021
$ws.Cells.Item(63, 7).Value = 0.05871433458377853
$ws.Cells.Item(63, 8).Value = 80.16338567370688
$ws.Cells.Item(64, 7).Value = 0.01800360959049803
$ws.Cells.Item(64, 8).Value = -55.57536105851422
$ws.Cells.Item(65, 7).Value = 0.05679254753738667
$ws.Cells.Item(65, 8).Value = 1.303248931741971
$ws.Cells.Item(66, 7).Value = 0.1013912121159845
$ws.Cells.Item(66, 8).Value = 8.377085710302403
$ws.Cells.Item(67, 7).Value = 0.08861652865372539
$ws.Cells.Item(67, 8).Value = -23.24028135950991
$ws.Cells.Item(68, 7).Value = -0.01483484251583027
$ws.Cells.Item(68, 8).Value = 57.43289247510354
$ws.Cells.Item(69, 7).Value = -0.02108640432576117
$ws.Cells.Item(69, 8).Value = 0.6382503037910543
$ws.Cells.Item(70, 7).Value = 0.07418125034718016
$ws.Cells.Item(70, 8).Value = -19.922013277562
$ws.Cells.Item(71, 7).Value = 0.0943058246545119
$ws.Cells.Item(71, 8).Value = 3.395536924504248
$ws.Cells.Item(72, 7).Value = -0.06704701991180746
$ws.Cells.Item(72, 8).Value = -19.55383937339416
$ws.Cells.Item(73, 7).Value = -0.06196467059367657
$ws.Cells.Item(73, 8).Value = 15.99444826911324
$ws.Cells.Item(74, 7).Value = 0.08770351430618467
$ws.Cells.Item(74, 8).Value = -12.24980855467511
$ws.Cells.Item(75, 7).Value = 0.1071687980909943
$ws.Cells.Item(75, 8).Value = 10.02502271624788
$ws.Cells.Item(76, 7).Value = 0.01168001437719657
$ws.Cells.Item(76, 8).Value = -54.32248314491348
$ws.Cells.Item(77, 7).Value = 0.03419671685587713
$ws.Cells.Item(77, 8).Value = 142.37450951045
$ws.Cells.Item(78, 7).Value = 0.07158858032586138
$ws.Cells.Item(78, 8).Value = 11.37516573386641
$ws.Cells.Item(79, 7).Value = 0.07029807107942133
$ws.Cells.Item(79, 8).Value = -8.363315480634673
$ws.Cells.Item(80, 7).Value = -0.1703849009931187
$ws.Cells.Item(80, 8).Value = -2.884212753485146
$ws.Cells.Item(81, 7).Value = -0.1473686933777397
$ws.Cells.Item(81, 8).Value = 29.855166689901
$ws.Cells.Item(82, 7).Value = 0.1488791644882012
$ws.Cells.Item(82, 8).Value = 29.80403566482165
$ws.Cells.Item(83, 7).Value = 0.1953565089928853
$ws.Cells.Item(83, 8).Value = 9.761816218283395
$ws.Cells.Item(84, 7).Value = 0.0756315108316095
$ws.Cells.Item(84, 8).Value = 217.2775582501963
$ws.Cells.Item(85, 7).Value = 0.0650683303849946
$ws.Cells.Item(85, 8).Value = 5.671679826455108
